# Weekly data refresh: a new price observation is inserted as row 112
# (the existing rows 112-158 shift down to 113-159, and the former last
# row 158 becomes the new last row 159 - dimension grows from R158 to R159).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 112, pushing everything below it down.
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new weekly observation.
$ws.Range('A112').Value = 9
$ws.Range('B112').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C112').Value = 'Metropolitana'
$ws.Range('D112').Value = 44510
$ws.Range('E112').Value = 13
$ws.Range('F112').Value = 100112030
$ws.Range('G112').Value = 'Poroto granado'
$ws.Range('H112').Value = 'Sin especificar'
$ws.Range('I112').Value = 'Primera'
$ws.Range('J112').Value = 16
$ws.Range('K112').Value = 34000
$ws.Range('L112').Value = 36000
$ws.Range('M112').Value = 35000
$ws.Range('N112').Value = '$/malla 25 kilos'
$ws.Range('O112').Value = 'Perú'
$ws.Range('P112').Value = 1400
$ws.Range('Q112').Value = 25
$ws.Range('R112').Value = 'Hortaliza'
